$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 541; this shifts the existing rows
# 541-636 down to 542-637 (matching the dimension change A1:R636 -> A1:R637).
$ws.Rows.Item(541).Insert()

# Populate the newly inserted row 541 with the new data point.
$ws.Cells.Item(541, 1).Value = 4
$ws.Cells.Item(541, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(541, 3).Value = "Los Lagos"
$ws.Cells.Item(541, 4).Value = 45180
$ws.Cells.Item(541, 5).Value = 10
$ws.Cells.Item(541, 6).Value = 100114013
$ws.Cells.Item(541, 7).Value = "Zanahoria"
$ws.Cells.Item(541, 8).Value = "Sin especificar"
$ws.Cells.Item(541, 9).Value = "Primera"
$ws.Cells.Item(541, 10).Value = 150
$ws.Cells.Item(541, 11).Value = 7000
$ws.Cells.Item(541, 12).Value = 7000
$ws.Cells.Item(541, 13).Value = 7000
$ws.Cells.Item(541, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(541, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(541, 16).Value = 350
$ws.Cells.Item(541, 17).Value = 20
$ws.Cells.Item(541, 18).Value = "Hortaliza"
